$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (PersonCrmRt), shifting existing
# columns C:F to D:G. This mirrors Excel's "Insert Column" behavior which
# shifts cells and carries the left neighbor's column formatting.
$ws.Columns("C").Insert()

# New column header
$ws.Range("C1").Value = "PopulationCen"

# New column data (PopulationCen values)
$values = @(585436, 595410, 604285, 609970, 620647, 631539, 643136, 648630, 643115)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

# New column inherits column B's (typed, non-bestFit) width; columns D:G
# keep the exact widths they had before the insert (handled automatically
# by Columns.Insert shifting the existing <col> entries).
$ws.Columns("C").ColumnWidth = 12.166666666666666

# Update selection to match final state
$ws.Range("E21").Select()
